# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on Hoja1!A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$note = @"
Conversión del día 💰
✅ Dólar paralelo: 68

Binance
✅ 1000 Bs = 12.35 = 49790.12 pesos
✅ 49790.12 pesos = 12.31 = 977.46 Bs

Promedio competencia
✅ Tasa pesos: 20
✅ Tasa Bs: 20
✅ % Ganancia: 20%
"@
$wsHoja1.Range("A1").Value = $note

# --- Update the rate cells on the "tasas" sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 81
$wsTasas.Range("O10").Value = 4033
$wsTasas.Range("N12").Value = 4045
$wsTasas.Range("O12").Value = 79.41
